$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows appended at the bottom of the sheet (rows 167-170),
# copying the formatting (bold/border for col A, date number-format for col E)
# from the last existing data row (166).
$ws.Range("A166").Copy()
$ws.Range("A167:A170").PasteSpecial(-4122)
$ws.Range("E166").Copy()
$ws.Range("E167:E170").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=167; A=165; B=6992690; E=45381.33333333334; F="Uthai Thani FC"; G="Prachuap FC";
       K=2.3; L=3.2; M=2.7; N=2.4; O=3.2; P=2.55; Q=0; R=1.8; S=2; T=2.75; U=1.9; V=1.9; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=168; A=166; B=6992689; E=45381.35416666666; F="Lamphun Warrior FC"; G="Port FC";
       K=3.25; L=3.5; M=1.909; N=3.5; O=3.4; P=1.85; Q=0.5; R=1.875; S=1.925; T=3; U=2; V=1.8; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=169; A=167; B=6992691; E=45381.375; F="Sukhothai FC"; G="Ratchaburi FC";
       K=2.75; L=3.25; M=2.25; N=3; O=3.25; P=2.1; Q=0.25; R=1.925; S=1.875; T=2.75; U=1.875; V=1.925; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=170; A=168; B=6992694; E=45381.41666666666; F="Trat FC"; G="BG Pathum United";
       K=4; L=4; M=1.615; N=3.75; O=3.8; P=1.666; Q=0.75; R=1.9; S=1.9; T=3.25; U=1.975; V=1.825; W=0; X=0; Y=0; Z=0; AA=0 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = "Thailand Premier League"
    $ws.Cells.Item($row, 4).Value = "Thailand Premier League"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
    $ws.Cells.Item($row, 23).Value = $r.W
    $ws.Cells.Item($row, 24).Value = $r.X
    $ws.Cells.Item($row, 25).Value = $r.Y
    $ws.Cells.Item($row, 26).Value = $r.Z
    $ws.Cells.Item($row, 27).Value = $r.AA
}
